# Workbook: Hortaliza, Vega Central Mapocho de Santiago - Arveja Verde
# A new weekly record is inserted as row 149, pushing all following rows
# down by one (old row 149 -> new row 150, ..., old row 190 -> new row 191).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149; this shifts rows 149..190 down to 150..191
# and preserves their existing formatting/values.
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new record's data.
$ws.Range("A149").Value = 9
$ws.Range("B149").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C149").Value = "Metropolitana"
$ws.Range("D149").Value = 45204
$ws.Range("E149").Value = 13
$ws.Range("F149").Value = 100112022
$ws.Range("G149").Value = "Arveja Verde"
$ws.Range("H149").Value = "Sin especificar"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 52
$ws.Range("K149").Value = 28000
$ws.Range("L149").Value = 29000
$ws.Range("M149").Value = 28558
$ws.Range("N149").Value = "$/malla 25 kilos"
$ws.Range("O149").Value = "Provincia de Limarí"
$ws.Range("P149").Value = 1142
$ws.Range("Q149").Value = 25
$ws.Range("R149").Value = "Hortaliza"
